# Apply commit "set index to Room ID in final df":
#   - readme sheet (Table1): reorder columns to index, Date, Author, sheet_name, JobNo
#   - each "Results, Air Speed *" sheet (Table4..Table12): swap Room Name / Room ID
#     columns so Room ID comes first, Room Name second.
#
# Data rows are moved with Range.Copy (through a scratch worksheet) so that
# text-typed cells (e.g. the "20220308" date string) keep their original
# string type instead of being re-interpreted as numbers. Header cells are
# assigned directly with .Value so the ListObject header text - and the
# table's column-name metadata - stay in sync.

$wb = $excel.ActiveWorkbook

# Scratch worksheet used purely as a round-trip buffer for Copy operations;
# removed again before the script finishes.
$helper = $wb.Worksheets.Add()

# ---------------------------------------------------------------------
# 1) "readme" sheet / Table1: index, Author, JobNo, sheet_name, Date
#                          -> index, Date, Author, sheet_name, JobNo
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("readme")
$lastRow = $ws.ListObjects.Item(1).Range.Rows.Count  # header + data rows

# Header row text (direct assignment keeps the ListObject / table column
# names synced, and none of these headers look numeric so the type stays text)
$ws.Cells.Item(1, 2).Value = "Date"
$ws.Cells.Item(1, 3).Value = "Author"
$ws.Cells.Item(1, 5).Value = "JobNo"
# column 4 (sheet_name) is unchanged

# Data rows (2..lastRow): move Author/JobNo/Date values around via the
# helper sheet so their shared-string / text typing is preserved.
$ws.Range("B2:B" + $lastRow).Copy($helper.Range("B2:B" + $lastRow))  # stash Author
$ws.Range("C2:C" + $lastRow).Copy($helper.Range("C2:C" + $lastRow))  # stash JobNo
$ws.Range("E2:E" + $lastRow).Copy($helper.Range("E2:E" + $lastRow))  # stash Date

$helper.Range("E2:E" + $lastRow).Copy($ws.Range("B2:B" + $lastRow))  # B = Date
$helper.Range("B2:B" + $lastRow).Copy($ws.Range("C2:C" + $lastRow))  # C = Author
$helper.Range("C2:C" + $lastRow).Copy($ws.Range("E2:E" + $lastRow))  # E = JobNo

$helper.Range("B2:E" + $lastRow).ClearContents()

# ---------------------------------------------------------------------
# 2) Results sheets / Table4..Table12: Room Name, Room ID, ...
#                                    -> Room ID, Room Name, ...
# ---------------------------------------------------------------------
$resultSheets = @(
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

foreach ($name in $resultSheets) {
    $rws = $wb.Worksheets.Item($name)
    $rLastRow = $rws.ListObjects.Item(1).Range.Rows.Count

    # Header row
    $rws.Cells.Item(1, 1).Value = "Room ID"
    $rws.Cells.Item(1, 2).Value = "Room Name"

    # Data rows: swap column A (Room Name) and column B (Room ID)
    $rws.Range("A2:A" + $rLastRow).Copy($helper.Range("A2:A" + $rLastRow))
    $rws.Range("B2:B" + $rLastRow).Copy($rws.Range("A2:A" + $rLastRow))
    $helper.Range("A2:A" + $rLastRow).Copy($rws.Range("B2:B" + $rLastRow))

    $helper.Range("A2:A" + $rLastRow).ClearContents()
}

$helper.Delete()
